# Update the "Next Week" (column F) entries for each team member's row so the
# spreadsheet matches the revised sprint-sheet document (commit: "update PDF,
# and update excel to match document").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("April 1")

$ws.Range("F9").Value  = "Client UI and server-side testing + project documentation (This doc, and README)"
$ws.Range("F10").Value = "Server – NodeJS Server & Client communication via strings & accept multiple clients"
$ws.Range("F11").Value = "Client & Server – Set up client-side w/ React & establish connection w/ NodeJS "
$ws.Range("F12").Value = "Server – Front end keyboard input and classes for server side"

# Reflect the author's final cursor position / selection in the saved view state.
$ws.Range("F12").Select()
